$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row 9 (continuation of the "Pattern Gen" task entries),
# copying the formatting down from row 8 so the new row matches the
# existing table styling (date/time columns, task columns, etc.).
$ws.Rows.Item(9).Insert(-4121, 0)   # xlShiftDown, xlFormatFromLeftOrAbove

$ws.Range("A9").Value = "29.3.2020"
$ws.Range("B9").Value = 0.82291666666666663
$ws.Range("C9").Value = 0.83333333333333337
$ws.Range("D9").Formula = "=C9-B9"
$ws.Range("E9").Value = "Pattern Gen 2"
$ws.Range("F9").Value = "Create necessary files"

# New, still-empty row 10 right below, matching column A's formatting
# (this is where the cursor ends up after entering the row above).
$ws.Range("A9").Copy() | Out-Null
$ws.Range("A10").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("A10").ClearContents() | Out-Null

# Move the active cell/selection to A10, matching the new end-of-data position
$ws.Range("A10").Select() | Out-Null
